# Updates crypto price/volume(1h) figures from the latest GitHub Actions scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'321.78"
$ws.Range("E2").Value = "'7.63%"
$ws.Range("D3").Value = "'48.69"
$ws.Range("E3").Value = "'15.38%"
$ws.Range("D4").Value = "'5.269"
$ws.Range("E4").Value = "'4.97%"
$ws.Range("D5").Value = "'0.08115"
$ws.Range("E5").Value = "'7.80%"
$ws.Range("D6").Value = "'4.580"
$ws.Range("E6").Value = "'5.03%"
$ws.Range("D7").Value = "'1.646"
$ws.Range("E7").Value = "'2.94%"
$ws.Range("E8").Value = "'28.57%"
$ws.Range("D9").Value = "'0.1297"
$ws.Range("E9").Value = "'9.06%"
$ws.Range("E10").Value = "'5.75%"
$ws.Range("D11").Value = "'0.09491"
$ws.Range("E11").Value = "'5.12%"
$ws.Range("D12").Value = "'0.04622"
$ws.Range("E12").Value = "'11.19%"
$ws.Range("D13").Value = "'0.1050"
$ws.Range("E13").Value = "'0.40%"
$ws.Range("D14").Value = "'0.001322"
$ws.Range("E14").Value = "'3.22%"
$ws.Range("D15").Value = "'0.04160"
$ws.Range("E15").Value = "'1.55%"
$ws.Range("D16").Value = "'0.005934"
$ws.Range("E16").Value = "'0.63%"
$ws.Range("D17").Value = "'3.343"
$ws.Range("E17").Value = "'0.12%"
$ws.Range("E18").Value = "'1.68%"
$ws.Range("E19").Value = "'2.04%"
$ws.Range("D20").Value = "'8.079"
$ws.Range("E20").Value = "'-2.57%"
$ws.Range("D21").Value = "'0.1400"
$ws.Range("E21").Value = "'0.71%"
$ws.Range("D22").Value = "'0.3124"
$ws.Range("E22").Value = "'0.73%"
$ws.Range("D23").Value = "'0.001306"
$ws.Range("E23").Value = "'3.18%"
$ws.Range("D24").Value = "'0.004254"
$ws.Range("E24").Value = "'9.22%"
$ws.Range("D25").Value = "'0.0001350"
$ws.Range("E25").Value = "'3.84%"
$ws.Range("D26").Value = "'0.0003539"
$ws.Range("E26").Value = "'-4.97%"
$ws.Range("D38").Value = "'0.02719"
$ws.Range("E38").Value = "'12.89%"
$ws.Range("D39").Value = "'0.05772"
$ws.Range("E39").Value = "'10.49%"
$ws.Range("D40").Value = "'0.006300"
$ws.Range("E40").Value = "'-0.01%"
$ws.Range("D41").Value = "'0.007681"
$ws.Range("E41").Value = "'-1.04%"
$ws.Range("D42").Value = "'0.1442"
$ws.Range("E42").Value = "'8.62%"
$ws.Range("D43").Value = "'0.007694"
$ws.Range("E43").Value = "'4.00%"
$ws.Range("E44").Value = "'13.65%"
$ws.Range("E45").Value = "'6.62%"
$ws.Range("D46").Value = "'0.00007008"
$ws.Range("E46").Value = "'12.51%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.01%"
$ws.Range("D48").Value = "'0.06255"
$ws.Range("E48").Value = "'37.84%"
$ws.Range("D49").Value = "'0.004000"
$ws.Range("E49").Value = "'-4.78%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.01%"
